# Append: 2025-10-11 18:20 JST
# Refreshes the "取得日時" capture timestamp for every existing row on the
# "ランサーズ" sheet and inserts 3 newly scraped Lancers job postings, pushing
# the previously-last rows further down the list (sheet grows from H10 to H13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @{ Row = 2; A = "2025-10-11 18:20:21"; B = "【AI開発者募集】多機能転売ツールの構築をお願いします!"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5407785"; G = 420; H = "🔥AI,Ai ◆ツール,開発" },
    @{ Row = 3; A = "2025-10-11 18:20:21"; B = "【急募】紙の伝票をWEBシステムへ自動データ入力開発【AI使用可能】"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5411519"; G = 383; H = "🔥AI,Ai ◆開発" },
    @{ Row = 4; A = "2025-10-11 18:20:21"; B = "急募 PR Zoom/Meet×TLDV×ChatGPT×Notion×Slack 議事録ワークフロー構築依頼"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5410688"; G = 323; H = "🔥GPT,ChatGPT" },
    @{ Row = 5; A = "2025-10-11 18:20:21"; B = "【急募】仕事の予約システム構築をお手伝いください!【AI使用可能】"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5411525"; G = 318; H = "🔥AI,Ai" },
    @{ Row = 6; A = "2025-10-11 18:20:21"; B = "【急募】配送状況を自動取得するAPI開発者募集!"; C = "システム開発"; D = "20,000 円 ~ 50,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5411268"; G = 238; H = "🔥API ◆開発" },
    @{ Row = 7; A = "2025-10-11 18:20:21"; B = "【自動売買】Excelと楽天RSSを活用したシステム開発依頼"; C = "システム開発"; D = "5,000 円 ~ 10,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5411684"; G = 110; H = "◆開発,システム開発" },
    @{ Row = 8; A = "2025-10-11 18:20:21"; B = "【急募】Webアプリ開発エンジニア募集!フルリモート可"; C = "システム開発"; D = "200,000 円 ~ 300,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5411585"; G = 93; H = "◆開発 ◇アプリ" },
    @{ Row = 9; A = "2025-10-11 18:20:21"; B = "急募バックエンドエンジニア マッチングサイトの開発"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5404059"; G = 93; H = "◆開発 ◇サイト" },
    @{ Row = 10; A = "2025-10-11 18:20:21"; B = "スプレッドシートをもとにした顧客・売上管理アプリのグライド化(Glide/無料版)"; C = "システム開発"; D = "5,000 円 ~ 10,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5411304"; G = 55; H = "◇アプリ" },
    @{ Row = 11; A = "2025-10-11 18:20:21"; B = "【急募】教育系のWEBサイトの作成"; C = "システム開発"; D = "20,000 円 ~ 50,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5411679"; G = 33; H = "◇サイト" },
    @{ Row = 12; A = "2025-10-11 18:20:21"; B = "【急募】時間単位で入札できるシステム構築の依頼"; C = "システム開発"; D = "100,000 円 ~ 200,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5411365"; G = 33; H = $null },
    @{ Row = 13; A = "2025-10-11 18:20:21"; B = "【フォーム制作】物件見積り査定フォーム制作の依頼"; C = "システム開発"; D = "20,000 円 ~ 50,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5411435"; G = 13; H = $null }
)

# Write all row values (A:H) in one pass. Rows 12 and 13 intentionally have no
# skill-summary text (column H), matching the source rows that were pushed down.
foreach ($row in $rows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    if ($null -eq $row.H) {
        $ws.Cells.Item($r, 8).ClearContents()
    } else {
        $ws.Cells.Item($r, 8).Value = $row.H
    }
}

# Rebuild the F-column hyperlinks (URL target must track the row it now belongs
# to) and re-apply the built-in "Hyperlink" cell style used throughout the sheet.
$ws.Hyperlinks.Delete()
foreach ($row in $rows) {
    $cell = $ws.Cells.Item($row.Row, 6)
    $ws.Hyperlinks.Add($cell, $row.F)
    $cell.Style = "Hyperlink"
}

